# Add a new "Save" column (H) to the s_vals sheet, matching the header
# styling used by the existing columns, and fill in the values for the
# two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" with the same style as the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells H2/H3: numeric 0, no special style (same as rest of data columns)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
